# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns
# for the coin rows on Sheet1, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Ref = "D2"; Value = "91.380.50" }
    @{ Ref = "E2"; Value = "  -0.57%  " }
    @{ Ref = "D3"; Value = "3.106.22" }
    @{ Ref = "E3"; Value = "  -0.12%  " }
    @{ Ref = "E4"; Value = "  +0.12%  " }
    @{ Ref = "D5"; Value = "242.57" }
    @{ Ref = "E5"; Value = "  -0.29%  " }
    @{ Ref = "D6"; Value = "614.53" }
    @{ Ref = "E6"; Value = "  -1.22%  " }
    @{ Ref = "E7"; Value = "  -3.24%  " }
    @{ Ref = "D8"; Value = "0.383" }
    @{ Ref = "E8"; Value = "  +2.72%  " }
    @{ Ref = "E9"; Value = "  +0.04%  " }
    @{ Ref = "D10"; Value = "3.107.34" }
    @{ Ref = "E10"; Value = "  +13.90%  " }
    @{ Ref = "E11"; Value = "  -2.91%  " }
    @{ Ref = "E12"; Value = "  +0.38%  " }
    @{ Ref = "D13"; Value = "0.0000248" }
    @{ Ref = "E13"; Value = "  -2.06%  " }
    @{ Ref = "D14"; Value = "5.59" }
    @{ Ref = "E14"; Value = "  +2.30%  " }
    @{ Ref = "D15"; Value = "34.35" }
    @{ Ref = "E15"; Value = "  -3.05%  " }
    @{ Ref = "D16"; Value = "91.423.73" }
    @{ Ref = "E16"; Value = "  -0.33%  " }
    @{ Ref = "D18"; Value = "3.102.63" }
    @{ Ref = "E18"; Value = "  +0.29%  " }
    @{ Ref = "D19"; Value = "3.69" }
    @{ Ref = "E19"; Value = "  -1.12%  " }
    @{ Ref = "D20"; Value = "14.70" }
    @{ Ref = "E20"; Value = "  +0.80%  " }
    @{ Ref = "E21"; Value = "  +0.02%  " }
    @{ Ref = "E22"; Value = "  -0.63%  " }
    @{ Ref = "D23"; Value = "9.25" }
    @{ Ref = "E23"; Value = "  +1.33%  " }
    @{ Ref = "E24"; Value = "  -7.36%  " }
    @{ Ref = "E25"; Value = "  -1.53%  " }
    @{ Ref = "D26"; Value = "88.29" }
    @{ Ref = "E26"; Value = "  -2.60%  " }
    @{ Ref = "D27"; Value = "11.63" }
    @{ Ref = "E27"; Value = "  -2.59%  " }
    @{ Ref = "D28"; Value = "3.287.06" }
    @{ Ref = "E28"; Value = "  +0.96%  " }
    @{ Ref = "E29"; Value = "  -0.14%  " }
    @{ Ref = "E30"; Value = "  +25.92%  " }
    @{ Ref = "D31"; Value = "0.231" }
    @{ Ref = "E31"; Value = "  -2.15%  " }
    @{ Ref = "E32"; Value = "  -9.67%  " }
    @{ Ref = "D33"; Value = "0.176" }
    @{ Ref = "E33"; Value = "  +3.46%  " }
    @{ Ref = "D34"; Value = "9.27" }
    @{ Ref = "E34"; Value = "  -1.11%  " }
    @{ Ref = "E35"; Value = "  -2.62%  " }
    @{ Ref = "D36"; Value = "7.63" }
    @{ Ref = "E36"; Value = "  -0.22%  " }
    @{ Ref = "D37"; Value = "26.08" }
    @{ Ref = "E37"; Value = "  -2.02%  " }
    @{ Ref = "D38"; Value = "4.02" }
    @{ Ref = "E38"; Value = "  -1.81%  " }
    @{ Ref = "E39"; Value = "  +0.88%  " }
    @{ Ref = "D40"; Value = "485.49" }
    @{ Ref = "E40"; Value = "  -0.85%  " }
    @{ Ref = "D41"; Value = "1.29" }
    @{ Ref = "E41"; Value = "  +0.55%  " }
    @{ Ref = "E42"; Value = "  +3.09%  " }
    @{ Ref = "D43"; Value = "3.39" }
    @{ Ref = "E43"; Value = "  -6.74%  " }
    @{ Ref = "E44"; Value = "  +0.11%  " }
    @{ Ref = "E45"; Value = "  +0.00%  " }
    @{ Ref = "D46"; Value = "159.44" }
    @{ Ref = "E46"; Value = "  +3.22%  " }
    @{ Ref = "D47"; Value = "0.695" }
    @{ Ref = "E47"; Value = "  +0.64%  " }
    @{ Ref = "D48"; Value = "1.88" }
    @{ Ref = "E48"; Value = "  -1.57%  " }
    @{ Ref = "E49"; Value = "  -0.82%  " }
    @{ Ref = "D50"; Value = "44.06" }
    @{ Ref = "E50"; Value = "  -1.25%  " }
    @{ Ref = "D51"; Value = "4.36" }
    @{ Ref = "E51"; Value = "  -4.34%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $value = $u.Value

    # The source data stores prices like "242.57" or "91.380.50" as plain
    # text (note some use "." as a thousands separator, producing values
    # such as "3.106.22" that are not valid numbers anyway). When a target
    # value would otherwise be auto-recognized by Excel as a pure number,
    # prefix it with an apostrophe so it is stored as text, just like the
    # neighboring non-numeric-looking price cells.
    if ($value -match '^[+-]?\d+(\.\d+)?$') {
        $cell.Value = "'" + $value
    } else {
        $cell.Value = $value
    }
}
